$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 0.86350148367952517
$ws.Range("F2").Value = 0.86486486486486491
$ws.Range("G2").Value = 0.86868686868686873

$ws.Range("E3").Value = 0.83381924198250723
$ws.Range("F3").Value = 0.81290322580645158
$ws.Range("G3").Value = 0.82935153583617749

[void]$ws.Columns("F").Select()

